$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1, matching style of B1 (bold, centered, bordered header)
$ws.Range("C1").Value = 2
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) # xlPasteFormats

# Data values for the new 't+3' column (C2:C12)
$values = @(
    -5.04614857794682,
    -1.194025718115943,
    -0.07161795042852842,
    -0.4111325302719243,
    0.0159267162195228,
    0.1002874912444511,
    0.1290666877551792,
    0.0276633633304105,
    0.03180697780879011,
    0.005376147938177376,
    0.03799774138790459
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
